$d = $word.ActiveDocument

$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center" /></w:pPr><w:r><w:rPr><w:b /><w:sz w:val="28" /></w:rPr><w:t>Dheeraj Chand</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center" /></w:pPr><w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2" /></w:pPr><w:r><w:t>PROFESSIONAL SUMMARY</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Senior Data Engineer with 21 years of expertise in geospatial data platforms, big data processing, and distributed systems architecture. Deep specialist in Apache Spark/Sedona for large-scale geospatial analytics, with fluency across ESRI, OSGeo, and SAFE FME technology stacks. Proven track record architecting production systems serving thousands of users, implementing PySpark pipelines processing billions of spatial records, and leading engineering teams. Expert in full-stack geospatial development from PostGIS database optimization to React-based mapping interfaces.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2" /></w:pPr><w:r><w:t>CORE COMPETENCIES</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Big Data &amp; Geospatial Processing: Apache Spark: PySpark, Spark SQL, Scala Spark, Sedona (geospatial), distributed processing • Geospatial Databases: PostGIS (advanced), Oracle Spatial, spatial indexing, query optimization • ETL/ELT: dbt, Informatica, CDAP, custom PySpark pipelines, data governance frameworks • Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Hadoop clusters, distributed computing • Streaming: Real-time data processing, Kafka integration, event-driven architectures</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>GIS Technology Stack: ESRI: ArcGIS Server, ArcGIS Pro, enterprise geodatabases, ModelBuilder, ArcPy scripting • OSGeo: QGIS, GRASS GIS, GDAL/OGR, GeoServer, spatial analysis workflows • SAFE FME: Data transformation, format conversion, spatial ETL, enterprise integration • Web Mapping: OpenLayers, Leaflet, MapBox, tile servers, WMS/WFS services • Spatial Analysis: Clustering algorithms, boundary estimation, network analysis, geostatistics</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Software Development &amp; Architecture: Python: Django/GeoDjango, Flask, Pandas, NumPy, SciKit-Learn, spatial libraries • JVM: Scala (Spark), Java (GeoTools, enterprise), Groovy scripting • Web Technologies: React, JavaScript, d3.js, RESTful APIs, microservices • Databases: PostgreSQL/PostGIS, Oracle, MySQL, MongoDB, spatial optimization • DevOps: Docker, Kubernetes, CI/CD (GitLab, GitHub), Airflow, Celery, nginx</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2" /></w:pPr><w:r><w:t>PROFESSIONAL EXPERIENCE</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>PARTNER &amp; SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 – Present</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Geospatial Data Platform Architecture and Big Data Engineering</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Architected and engineered production geospatial platforms serving thousands of analysts</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built enterprise-scale ETL pipelines using PySpark and Sedona processing billions of geospatial records with advanced spatial clustering algorithms</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed multi-tenant data warehouse integrating Census, electoral, and demographic data using PostGIS and Spark SQL optimization</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Implemented fraud detection systems processing multi-terabyte datasets with real-time spatial analysis capabilities</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Created parametric boundary estimation algorithms using PostGIS and GRASS without machine learning dependencies</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Led integration of ESRI ArcGIS Server, OSGeo tools (QGIS, GRASS), and SAFE FME for enterprise geospatial workflows</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 – 2014</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Geospatial Systems Architecture and Development</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Architected geospatial analysis frameworks and mapping applications for electoral research</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed custom visualization tools and interactive dashboards using JavaScript and OpenLayers</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Designed data processing pipelines for large-scale demographic and geographic datasets</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Implemented PostGIS spatial databases and optimized geospatial query performance</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 – 2012</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Data Platform Architecture and Engineering Leadership</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Led technical architecture and development of data-driven political technology platforms</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Designed and implemented scalable data platforms using Python, Django, and PostgreSQL</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built RESTful APIs and microservices architecture for campaign data integration</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Managed engineering teams and established development best practices and CI/CD pipelines</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>SENIOR TECHNICAL ANALYST - GSD&amp;M, Austin, TX | 2008 – 2010</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Campaign Technology and Data Engineering</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed voter targeting models and demographic analysis tools using Python and R</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built web applications and data visualization systems for campaign analytics</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Created data integration systems connecting multiple campaign data sources</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Implemented machine learning algorithms and statistical models for voter behavior prediction</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Political Technology and Data Systems</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Coordinated technical operations and data systems for political campaigns</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed data collection and analysis protocols for campaign research</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built custom applications and tools for voter engagement and campaign management</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Supported technical infrastructure and data processing for progressive political initiatives</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Political Technology Development</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed software solutions for political campaigns using PHP, JavaScript, and MySQL</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built web applications for voter engagement and campaign management</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Integrated third-party APIs and data sources for campaign tools</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Collaborated with political strategists to translate requirements into technical solutions</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Nonprofit Technology Integration and Development</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed data management systems and web applications for social justice organizations</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built custom applications for community engagement using PHP, MySQL, and web technologies</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Provided technical training and support to nonprofit staff</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Integrated technology solutions within organizational frameworks for advocacy work</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Political Research and Data Analysis Tools</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed data analysis tools for political polling and research using Python and R</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built statistical models and data visualization tools for research presentations</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Created automated reporting systems and data processing pipelines for survey analysis</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Supported senior researchers with technical analysis and data processing automation</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Field Operations Technology and Data Management</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed data collection and management systems for political field operations</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built databases and reporting tools for campaign field work and voter outreach</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Trained field staff on data collection protocols and quality control systems</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Analyzed field data using statistical methods to inform campaign strategy and research findings</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2" /></w:pPr><w:r><w:t>KEY ACHIEVEMENTS AND IMPACT</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3" /></w:pPr><w:r><w:t>Geospatial Platform Engineering</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Architected redistricting platform processing Census data for thousands of analysts with real-time PostGIS collaborative editing</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Built boundary estimation system using advanced PostGIS algorithms and incomplete data without machine learning requirements</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>• Developed geospatial simulation platform integrating multi-agent modeling with web interface</w:t></w:r></w:p>'

$d.Content.InsertXML($bodyXml)

$d.PageSetup.TopMargin = 72
$d.PageSetup.BottomMargin = 72
$d.PageSetup.LeftMargin = 90
$d.PageSetup.RightMargin = 90

Write-Output $d.Paragraphs.Count
